# Auto-generated Excel COM-interop script
# Applies a market-data refresh: updates price/profit columns (H:N)
# for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 117.78571
$ws.Range("I6").Value = 124.92308
$ws.Range("K6").Value = 374.76924
$ws.Range("M6").Value = -262.76924

# Row 8
$ws.Range("H8").Value = 76.666664
$ws.Range("I8").Value = 76.666664
$ws.Range("K8").Value = 229.999992
$ws.Range("M8").Value = -90.99999199999999

# Row 39
$ws.Range("H39").Value = 475.5
$ws.Range("J39").Value = 617.0833
$ws.Range("L39").Value = 1851.2499
$ws.Range("N39").Value = -2443.2499

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("N45").Value = 0

# Row 80
$ws.Range("H80").Value = 1287.4706
$ws.Range("I80").Value = 1160.4286
$ws.Range("J80").Value = 1376.4
$ws.Range("K80").Value = 3481.2858
$ws.Range("L80").Value = 4129.200000000001
$ws.Range("M80").Value = -2483.2858
$ws.Range("N80").Value = -6125.200000000001

# Row 83
$ws.Range("H83").Value = 1287.4706
$ws.Range("I83").Value = 1160.4286
$ws.Range("J83").Value = 1376.4
$ws.Range("K83").Value = 10443.8574
$ws.Range("L83").Value = 12387.6
$ws.Range("M83").Value = -5451.857399999999
$ws.Range("N83").Value = -22371.6

# Row 96
$ws.Range("H96").Value = 1686.125
$ws.Range("I96").Value = 583.3333
$ws.Range("J96").Value = 2347.8
$ws.Range("K96").Value = 1749.9999
$ws.Range("L96").Value = 7043.400000000001
$ws.Range("M96").Value = -376.9999
$ws.Range("N96").Value = -9789.400000000001

# Row 98
$ws.Range("H98").Value = 1616.8462
$ws.Range("I98").Value = 1418.25
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1418.25
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = 79.75
$ws.Range("N98").Value = -6996

# Row 112
$ws.Range("H112").Value = 2097.6
$ws.Range("I112").Value = 1497
$ws.Range("J112").Value = 2498
$ws.Range("K112").Value = 4491
$ws.Range("L112").Value = 7494
$ws.Range("M112").Value = -3383
$ws.Range("N112").Value = -9710

# Row 122
$ws.Range("H122").Value = 1616.8462
$ws.Range("I122").Value = 1418.25
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4254.75
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1804.75
$ws.Range("N122").Value = -16900

# Row 132
$ws.Range("H132").Value = 1134.35
$ws.Range("I132").Value = 1299.4375
$ws.Range("J132").Value = 474
$ws.Range("K132").Value = 3898.3125
$ws.Range("L132").Value = 1422
$ws.Range("M132").Value = -1368.3125
$ws.Range("N132").Value = -6482

# Row 137
$ws.Range("H137").Value = 2298.3333
$ws.Range("I137").Value = 1274.6875
$ws.Range("J137").Value = 4345.625
$ws.Range("K137").Value = 3824.0625
$ws.Range("L137").Value = 13036.875
$ws.Range("M137").Value = -1274.0625
$ws.Range("N137").Value = -18136.875

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2974.75
$ws.Range("J45").Value = 3571
$ws.Range("L45").Value = 3571
$ws.Range("N45").Value = -4325

# Row 74
$ws.Range("H74").Value = 3035.0908
$ws.Range("I74").Value = 3048.6
$ws.Range("J74").Value = 2900
$ws.Range("K74").Value = 3048.6
$ws.Range("L74").Value = 2900
$ws.Range("M74").Value = -2174.6
$ws.Range("N74").Value = -4648

# Row 77
$ws.Range("H77").Value = 3035.0908
$ws.Range("I77").Value = 3048.6
$ws.Range("J77").Value = 2900
$ws.Range("K77").Value = 15243
$ws.Range("L77").Value = 14500
$ws.Range("M77").Value = -10875
$ws.Range("N77").Value = -23236

# Row 132
$ws.Range("H132").Value = 1257.3125
$ws.Range("J132").Value = 2105
$ws.Range("L132").Value = 6315
$ws.Range("N132").Value = -11375

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1148.9166
$ws.Range("I107").Value = 1117
$ws.Range("K107").Value = 1117
$ws.Range("M107").Value = 803

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3184.0667
$ws.Range("I31").Value = 3090.2144
$ws.Range("J31").Value = 4498
$ws.Range("K31").Value = 3090.2144
$ws.Range("L31").Value = 4498
$ws.Range("M31").Value = -2795.2144
$ws.Range("N31").Value = -5088

# Row 34
$ws.Range("H34").Value = 3184.0667
$ws.Range("I34").Value = 3090.2144
$ws.Range("J34").Value = 4498
$ws.Range("K34").Value = 3090.2144
$ws.Range("L34").Value = 4498
$ws.Range("M34").Value = -2888.2144
$ws.Range("N34").Value = -4902

# Row 99
$ws.Range("H99").Value = 7754.091
$ws.Range("I99").Value = 5061.2
$ws.Range("J99").Value = 9998.166999999999
$ws.Range("K99").Value = 5061.2
$ws.Range("L99").Value = 9998.166999999999
$ws.Range("M99").Value = -3563.2
$ws.Range("N99").Value = -12994.167

# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -13900

# Row 126
$ws.Range("H126").Value = 7754.091
$ws.Range("I126").Value = 5061.2
$ws.Range("J126").Value = 9998.166999999999
$ws.Range("K126").Value = 15183.6
$ws.Range("L126").Value = 29994.501
$ws.Range("M126").Value = -12713.6
$ws.Range("N126").Value = -34934.501

# Row 132
$ws.Range("H132").Value = 8544.286
$ws.Range("I132").Value = 5962
$ws.Range("K132").Value = 17886
$ws.Range("M132").Value = -15356

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 1045.4166
$ws.Range("I12").Value = 111.666664
$ws.Range("J12").Value = 1356.6666
$ws.Range("K12").Value = 334.999992
$ws.Range("L12").Value = 4069.9998
$ws.Range("M12").Value = -161.999992
$ws.Range("N12").Value = -4415.9998

# Row 98
$ws.Range("H98").Value = 702.1429000000001
$ws.Range("I98").Value = 706.75
$ws.Range("J98").Value = 696
$ws.Range("K98").Value = 2120.25
$ws.Range("L98").Value = 2088
$ws.Range("M98").Value = -622.25
$ws.Range("N98").Value = -5084

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 7698
$ws.Range("I80").Value = 2396
$ws.Range("K80").Value = 2396
$ws.Range("M80").Value = -1398

# Row 83
$ws.Range("H83").Value = 7698
$ws.Range("I83").Value = 2396
$ws.Range("K83").Value = 11980
$ws.Range("M83").Value = -6988

# Row 102
$ws.Range("H102").Value = 1132.3636
$ws.Range("I102").Value = 945.875
$ws.Range("K102").Value = 945.875
$ws.Range("M102").Value = 676.125

# Row 122
$ws.Range("H122").Value = 5110.0527
$ws.Range("I122").Value = 5123.0586
$ws.Range("K122").Value = 15369.1758
$ws.Range("M122").Value = -12919.1758

# Row 132
$ws.Range("H132").Value = 2259.4614
$ws.Range("I132").Value = 1884.6086
$ws.Range("K132").Value = 5653.825800000001
$ws.Range("M132").Value = -3123.825800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2343.5
$ws.Range("I7").Value = 2050
$ws.Range("K7").Value = 2050
$ws.Range("M7").Value = -1938

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 0

# Row 22
$ws.Range("H22").Value = 3774.375
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 4199.2856
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 4199.2856
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -4789.2856

# Row 27
$ws.Range("H27").Value = 3774.375
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 4199.2856
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 4199.2856
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -4413.2856

# Row 40
$ws.Range("H40").Value = 11931.353
$ws.Range("I40").Value = 11855.6
$ws.Range("K40").Value = 11855.6
$ws.Range("M40").Value = -11719.6

# Row 68
$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

# Row 71
$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

# Row 82
$ws.Range("H82").Value = 7737.8
$ws.Range("I82").Value = 3845
$ws.Range("K82").Value = 3845
$ws.Range("M82").Value = -3484

# Row 85
$ws.Range("H85").Value = 7737.8
$ws.Range("I85").Value = 3845
$ws.Range("K85").Value = 3845
$ws.Range("M85").Value = -2597

# Row 122
$ws.Range("H122").Value = 8273.362999999999
$ws.Range("I122").Value = 8361.888999999999
$ws.Range("K122").Value = 25085.667
$ws.Range("M122").Value = -22635.667

# Row 126
$ws.Range("H126").Value = 2343.5
$ws.Range("I126").Value = 2050
$ws.Range("K126").Value = 6150
$ws.Range("M126").Value = -3680

# Row 136
$ws.Range("H136").Value = 5200
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 5333.3335
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 16000.0005
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -21100.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 1500
$ws.Range("J5").Value = 1500
$ws.Range("L5").Value = 1500
$ws.Range("N5").Value = -1724

# Row 62
$ws.Range("H62").Value = 7800.2
$ws.Range("J62").Value = 6501.5
$ws.Range("L62").Value = 6501.5
$ws.Range("N62").Value = -7749.5

# Row 65
$ws.Range("H65").Value = 7800.2
$ws.Range("J65").Value = 6501.5
$ws.Range("L65").Value = 32507.5
$ws.Range("N65").Value = -38747.5

# Row 122
$ws.Range("H122").Value = 1041.3334
$ws.Range("I122").Value = 1037.25
$ws.Range("J122").Value = 1049.5
$ws.Range("K122").Value = 3111.75
$ws.Range("L122").Value = 3148.5
$ws.Range("M122").Value = -661.75
$ws.Range("N122").Value = -8048.5

# Row 126
$ws.Range("H126").Value = 2623.5
$ws.Range("I126").Value = 1998.2
$ws.Range("K126").Value = 5994.6
$ws.Range("M126").Value = -3524.6
